# Auto-generated edit script: updates crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.328.04'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '1.609.73'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0614'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').Value = '1.833.47'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').Value = '1.620.79'
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('E15').Value = '  +1.12%  '
$ws.Range('D16').Value = '26.302.15'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.11'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.78%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -3.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  +2.05%  '
$ws.Range('E30').Value = '  +4.60%  '
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.19%  '
$ws.Range('E34').Value = '  +2.92%  '
$ws.Range('E35').Value = '  +1.31%  '
$ws.Range('D36').Value = '1.167.82'
$ws.Range('E36').Value = '  +5.60%  '
$ws.Range('E37').Value = '  +2.23%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.789'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.498'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.782'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.79%  '
$ws.Range('D44').Value = '1.749.69'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.22'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.62%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₇0994'
$ws.Range('E48').Value = '  -5.43%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0507'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.407'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('E51').Value = '  -0.06%  '

Write-Output "done"
